# Commit message: "Change inserting date from manually to grab it from
# website automatically".
#
# The sheet holds one day's worth of Industry-Indices rows (r=2..42).
# The "grabbed from website" values differ from the old manually-entered
# ones in two ways:
#   1. The CDate (col A) and JDate (col B) strings gain ISO-style dashes:
#        20190603  -> 2019-06-03
#        13980313  -> 1398-03-13
#   2. MarketValue/TransactionsVol/TransactionsValue (cols E, G, H) were
#      being over-scaled by 10^6 on manual entry; divide each by
#      1,000,000. For G, only values that were >= 1,000,000 get divided
#      (smaller volumes were already correctly scaled). When a division
#      leaves a fractional remainder, the site reports it as literal text
#      with 3 decimals and a trailing space (e.g. "172.604 ") instead of
#      a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 42

for ($r = $firstRow; $r -le $lastRow; $r++) {

    # --- Column A: CDate, e.g. 20190603 -> 2019-06-03 ---
    $aRaw = [string]$ws.Cells.Item($r, 1).Value2
    if ($aRaw.Length -eq 8) {
        $aNew = $aRaw.Substring(0, 4) + "-" + $aRaw.Substring(4, 2) + "-" + $aRaw.Substring(6, 2)
        $ws.Cells.Item($r, 1).Value = "'" + $aNew
    }

    # --- Column B: JDate, e.g. 13980313 -> 1398-03-13 ---
    $bRaw = [string]$ws.Cells.Item($r, 2).Value2
    if ($bRaw.Length -eq 8) {
        $bNew = $bRaw.Substring(0, 4) + "-" + $bRaw.Substring(4, 2) + "-" + $bRaw.Substring(6, 2)
        $ws.Cells.Item($r, 2).Value = "'" + $bNew
    }

    # --- Column E: MarketValue, always /1,000,000 ---
    $eRaw = $ws.Cells.Item($r, 5).Value2
    $eNew = [math]::Round($eRaw / 1000000, 6)
    $ws.Cells.Item($r, 5).Value = $eNew

    # --- Column G: TransactionsVol, /1,000,000 only when big enough ---
    $gRaw = $ws.Cells.Item($r, 7).Value2
    if ($gRaw -ge 1000000) {
        $gNew = [math]::Round($gRaw / 1000000, 6)
        if ($gNew -ne [math]::Floor($gNew)) {
            $gText = "{0:F3} " -f $gNew
            $ws.Cells.Item($r, 7).Value = "'" + $gText
        } else {
            $ws.Cells.Item($r, 7).Value = $gNew
        }
    }

    # --- Column H: TransactionsValue, always /1,000,000 ---
    $hRaw = $ws.Cells.Item($r, 8).Value2
    $hNew = [math]::Round($hRaw / 1000000, 6)
    if ($hNew -ne [math]::Floor($hNew)) {
        $hText = "{0:F3} " -f $hNew
        $ws.Cells.Item($r, 8).Value = "'" + $hText
    } else {
        $ws.Cells.Item($r, 8).Value = $hNew
    }
}
